$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Maharashtra", 12),
    @("Kerala", 2),
    @("Karnataka", 16),
    @("Tamil Nadu", 5),
    @("Andhra Pradesh", 3),
    @("West Bengal", 5),
    @("Delhi", 8),
    @("Odisha", 3),
    @("Rajasthan", 7),
    @("Gujarat", 2),
    @("Chhattisgarh", 15),
    @("Jammu and Kashmir", 5),
    @("Himachal Pradesh", 2),
    @("Goa", 2),
    @("Puducherry", 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Match the formatting of the newly-added rows (10-16) to the existing
# data rows (the original sheet only had data through row 9).
$ws.Range("A2:B2").Copy()
$ws.Range("A10:B16").PasteSpecial(-4122)  # xlPasteFormats
